# odczyt diagnostyki silnika na GUI
# Adds four new translation rows (49-52) to the "Translation" sheet for the
# new "Motor Diagn:" GUI text and its related single-use id placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 49: Motor Diagn: label
$ws.Range("B49").Value = "SingleUseId69"
$ws.Range("C49").Value = "Typography_00"
$ws.Range("D49").Value = "Center"
$ws.Range("E49").Value = "LTR"
$ws.Range("F49").Value = "Motor Diagn:"

# Row 50: value/value placeholder text reused from existing shared string
$ws.Range("B50").Value = "SingleUseId70"
$ws.Range("C50").Value = "Typography_00"
$ws.Range("D50").Value = "Center"
$ws.Range("E50").Value = "LTR"
$ws.Range("F50").Value = "<value>/<value>"

# Row 51: 65535 value (stored as text, matching the source workbook; the
# leading apostrophe forces text entry and Style="Normal" drops the
# auto-generated quote-prefix formatting so the cell keeps the default style)
$ws.Range("B51").Value = "SingleUseId71"
$ws.Range("C51").Value = "Typography_00"
$ws.Range("D51").Value = "Left"
$ws.Range("E51").Value = "LTR"
$ws.Range("F51").Value = "'65535"
$ws.Range("F51").Style = "Normal"

# Row 52: 65535 value (second occurrence, stored as text)
$ws.Range("B52").Value = "SingleUseId72"
$ws.Range("C52").Value = "Typography_00"
$ws.Range("D52").Value = "Left"
$ws.Range("E52").Value = "LTR"
$ws.Range("F52").Value = "'65535"
$ws.Range("F52").Style = "Normal"
